$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / account holder details
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 12.03.2025"

# Row 6
$ws.Range("B6").Value = "15.03."
$ws.Range("C6").Value = "16.03."
$ws.Range("D6").Value = "ZALANDO MKTPLC EU ZPAXCJ"
$ws.Range("E6").Value = "196,25-"

# Row 7
$ws.Range("B7").Value = "18.03."
$ws.Range("C7").Value = "19.03."
$ws.Range("D7").Value = "PAYPAL ZVOKYX"
$ws.Range("E7").Value = "74,64-"

# Row 8
$ws.Range("B8").Value = "20.03."
$ws.Range("C8").Value = "21.03."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "25,08-"

# Row 9 - newly populated (was empty)
$ws.Range("B9").Value = "22.03."
$ws.Range("C9").Value = "23.03."
$ws.Range("D9").Value = "EBAY MKTPLC EU GUIOKY"
$ws.Range("E9").Value = "202,14-"
$ws.Range("E9").Style = $ws.Range("E8").Style

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 24.03.2025"
$ws.Range("E12").Value = "498,11-"

# Next settlement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 29.03.2025"
